$wb = $excel.ActiveWorkbook

# Sheet "y_fitted_on_begin_2016" (sheet1) - update column B (y_value)
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Cells.Item(2, 2).Value = 16.2485988578499
$ws1.Cells.Item(3, 2).Value = 16.55658146164153
$ws1.Cells.Item(4, 2).Value = 16.53242161070152
$ws1.Cells.Item(5, 2).Value = 16.18166463077969
$ws1.Cells.Item(6, 2).Value = 16.35898625681103
$ws1.Cells.Item(7, 2).Value = 16.66453429402533
$ws1.Cells.Item(8, 2).Value = 16.38503688173556
$ws1.Cells.Item(9, 2).Value = 15.95575588151746
$ws1.Cells.Item(10, 2).Value = 16.08032003856645
$ws1.Cells.Item(11, 2).Value = 16.36155348947774
$ws1.Cells.Item(12, 2).Value = 16.79765429799629
$ws1.Cells.Item(13, 2).Value = 16.61200119478694
$ws1.Cells.Item(14, 2).Value = 15.54003571945567
$ws1.Cells.Item(15, 2).Value = 15.40563892040851
$ws1.Cells.Item(16, 2).Value = 16.41824720020411
$ws1.Cells.Item(17, 2).Value = 17.24845705007934
$ws1.Cells.Item(18, 2).Value = 17.42940042975582
$ws1.Cells.Item(19, 2).Value = 16.87224372638843
$ws1.Cells.Item(20, 2).Value = 15.85825509430436
$ws1.Cells.Item(21, 2).Value = 15.38184146475877
$ws1.Cells.Item(22, 2).Value = 15.63080864992266
$ws1.Cells.Item(23, 2).Value = 16.41002665604104
$ws1.Cells.Item(24, 2).Value = 16.10240848857352
$ws1.Cells.Item(25, 2).Value = 15.33062925660056
$ws1.Cells.Item(26, 2).Value = 15.61797843209589
$ws1.Cells.Item(27, 2).Value = 16.51831064661114
$ws1.Cells.Item(28, 2).Value = 16.53529536440785
$ws1.Cells.Item(29, 2).Value = 16.64097001610006
$ws1.Cells.Item(30, 2).Value = 16.78521562037126
$ws1.Cells.Item(31, 2).Value = 16.7148550440743
$ws1.Cells.Item(32, 2).Value = 16.68220913675773
$ws1.Cells.Item(33, 2).Value = 16.68342571737815
$ws1.Cells.Item(34, 2).Value = 16.09022478733687
$ws1.Cells.Item(35, 2).Value = 15.50938226738423
$ws1.Cells.Item(36, 2).Value = 15.50822895520085
$ws1.Cells.Item(37, 2).Value = 15.72998366417209
$ws1.Cells.Item(38, 2).Value = 16.03359202125094

# Sheet "y_fitted_on_begin_2021" (sheet3) - update column B (y_value)
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Cells.Item(2, 2).Value = 16.08944552471046
$ws3.Cells.Item(3, 2).Value = 16.66115993175783
$ws3.Cells.Item(4, 2).Value = 16.65192218058955
$ws3.Cells.Item(5, 2).Value = 16.19114867105358
$ws3.Cells.Item(6, 2).Value = 16.56508355755991
$ws3.Cells.Item(7, 2).Value = 16.54671281176488
$ws3.Cells.Item(8, 2).Value = 16.22560625568105
$ws3.Cells.Item(9, 2).Value = 16.05853029884547
$ws3.Cells.Item(10, 2).Value = 16.07345140303136
$ws3.Cells.Item(11, 2).Value = 15.9677962956466
$ws3.Cells.Item(12, 2).Value = 16.79964234344508
$ws3.Cells.Item(13, 2).Value = 16.35264215995071
$ws3.Cells.Item(14, 2).Value = 15.4682715156059
$ws3.Cells.Item(15, 2).Value = 15.7381207904387
$ws3.Cells.Item(16, 2).Value = 15.84605607145817
$ws3.Cells.Item(17, 2).Value = 16.69030076951861
$ws3.Cells.Item(18, 2).Value = 17.63924525182826
$ws3.Cells.Item(19, 2).Value = 17.15709447613781
$ws3.Cells.Item(20, 2).Value = 16.07823724989035
$ws3.Cells.Item(21, 2).Value = 15.41739337159053
$ws3.Cells.Item(22, 2).Value = 15.01055249797808
$ws3.Cells.Item(23, 2).Value = 15.86134502669368
$ws3.Cells.Item(24, 2).Value = 15.17159538831786
$ws3.Cells.Item(25, 2).Value = 15.55956891243692
$ws3.Cells.Item(26, 2).Value = 15.57891455041512
$ws3.Cells.Item(27, 2).Value = 15.84551368789917
$ws3.Cells.Item(28, 2).Value = 15.97337558856278
$ws3.Cells.Item(29, 2).Value = 17.26085366863868
$ws3.Cells.Item(30, 2).Value = 16.6017093265964
$ws3.Cells.Item(31, 2).Value = 16.86670847951624
$ws3.Cells.Item(32, 2).Value = 16.82165251493591
$ws3.Cells.Item(33, 2).Value = 16.66047115855641
$ws3.Cells.Item(34, 2).Value = 15.76940606821367
$ws3.Cells.Item(35, 2).Value = 15.57671541462953
$ws3.Cells.Item(36, 2).Value = 15.08369768791743
$ws3.Cells.Item(37, 2).Value = 15.02624484843359
$ws3.Cells.Item(38, 2).Value = 15.6077229691664
$ws3.Cells.Item(39, 2).Value = 15.16485051042688
$ws3.Cells.Item(40, 2).Value = 15.14146405613493
$ws3.Cells.Item(41, 2).Value = 15.07990969697848
$ws3.Cells.Item(42, 2).Value = 13.79863709251337
$ws3.Cells.Item(43, 2).Value = 14.31740365201596
